# normalize geometry before feeding to network
#
# Adds a new model run ("pointsemantic" with geometry bucket "[2:3]") to the
# semantic3d sheet: finishes the previously-incomplete row 13 ("[1:2]"
# bucket) and appends a new row 14 for the "[2:3]" bucket.

$wb = $excel.ActiveWorkbook

# --- sheet "semantic3d" (first sheet) ---------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()

# Row 13 ("pointsemantic" / "[1:2]") was missing its metric columns - fill
# them in now.
$ws1.Range("I13").Value = 40
$ws1.Range("J13").Value = 0.63315699999999997
$ws1.Range("K13").Value = 0.88558999999999999
$ws1.Range("L13").Value = 0.87949500000000003
$ws1.Range("M13").Value = 0.65441899999999997
$ws1.Range("N13").Value = 0.88510999999999995
$ws1.Range("O13").Value = 0.37359900000000001
$ws1.Range("P13").Value = 0.89678000000000002
$ws1.Range("Q13").Value = 0.30763000000000001
$ws1.Range("R13").Value = 0.534501
$ws1.Range("S13").Value = 0.53371900000000005

# New row 14: same model, next geometry bucket "[2:3]".
$ws1.Range("A14").Value = "pointsemantic"
$ws1.Range("B14").Value = 8192
$ws1.Range("C14").Value = 16
$ws1.Range("D14").Value = 1
$ws1.Range("E14").Value = 1
$ws1.Range("F14").Value = "[2:3]"
$ws1.Range("G14").Value = 10
$ws1.Range("H14").Value = 10
$ws1.Range("I14").Value = 15
$ws1.Range("J14").Value = 0.487541
$ws1.Range("K14").Value = 0.813558
$ws1.Range("L14").Value = 0.82316500000000004
$ws1.Range("M14").Value = 0.61223399999999994
$ws1.Range("N14").Value = 0.80680600000000002
$ws1.Range("O14").Value = 0.23685800000000001
$ws1.Range("P14").Value = 0.74796099999999999
$ws1.Range("Q14").Value = 0.143812
$ws1.Range("R14").Value = 0.24742400000000001
$ws1.Range("S14").Value = 0.28206599999999998

$ws1.Range("P20").Select()

# --- sheet "npm3d" (second sheet) --------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("T1").Select()

# --- sheet "common_class" (third sheet) --------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()
$ws3.Range("J18").Select()

# Leave the originally active sheet selected again.
$ws1.Activate()
